$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Problem 2) grades, entered top-to-bottom so that new
# shared strings are created in the same order as the source workbook. ---

$ws.Range("C2").Value = "10 (Fabios Solution?)"
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = "3 (general idea, but code does not work)"
$ws.Range("C7").Value = "1 (the code did not work. The code is complex, I did not follow the code)"
# C8 is intentionally left blank but highlighted (missing submission/grade)
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = "8 (general idea correct, but wrong check in the end)"
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = "9 (input from command line)"
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = "0 (no code)"
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = "3 (seems that you did not understand the concept of jolly)"
$ws.Range("C18").Value = "9 (need to convert inputs to integers)"
$ws.Range("C19").Value = "9 (any sequence with a single number is jolly!)"
$ws.Range("E1").Value = "Obs"
$ws.Range("E18").Value = "All your python files have merge conflicts, you should be careful when using more than one machine to work on the assignment."
$ws.Range("C20").Value = 10
$ws.Range("C21").Value = 10
$ws.Range("C22").Value = "3 (seems that you did not understand the concept of jolly)"
$ws.Range("C23").Value = "7 (write idea, but logic is wrong! And any sequence with a single number is jolly!)"
$ws.Range("C24").Value = 10
$ws.Range("C25").Value = "0 (no code)"
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = "0 (no code)"
$ws.Range("C28").Value = "9 (the first element of the sequence is the sequence size and fails for sequences with only one element)"
$ws.Range("C29").Value = "5 (seems that you did not understand well the concept of jolly, the code does not work in general)"
$ws.Range("C30").Value = 10
$ws.Range("C31").Value = 10
$ws.Range("C32").Value = "0 (no code)"
$ws.Range("C33").Value = 10
$ws.Range("C34").Value = "9 (Missing indentation in the nested for's)"
$ws.Range("C35").Value = 10

# Highlight the missing grade for C8 (yellow fill), matching the
# conditional-ish reminder styling used for ungraded cells.
$ws.Range("C8").Interior.Color = 65535

# --- Column widths for the new/updated columns ---
$ws.Columns.Item(3).ColumnWidth = 88.08984375
$ws.Columns.Item(5).ColumnWidth = 109.1796875

# --- Scroll/selection state left by the grader ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("C35").Select()
